# Weekly price update: a new week's record is inserted above the most
# recent entries for this market/product, pushing the existing rows down
# by one (row 552 -> 553, ..., row 607 -> 608) and extending the used
# range from A1:R607 to A1:R608.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 552; Excel shifts rows 552..607 down to 553..608
# and copies formatting (incl. the date style on column D) along with them.
$ws.Rows.Item(552).Insert()

# Populate the newly inserted row 552 with this week's record.
$ws.Cells.Item(552, 1).Value  = 4
$ws.Cells.Item(552, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(552, 3).Value  = "Los Lagos"
$ws.Cells.Item(552, 4).Value  = 45132
$ws.Cells.Item(552, 5).Value  = 10
$ws.Cells.Item(552, 6).Value  = 100114013
$ws.Cells.Item(552, 7).Value  = "Zanahoria"
$ws.Cells.Item(552, 8).Value  = "Sin especificar"
$ws.Cells.Item(552, 9).Value  = "Primera"
$ws.Cells.Item(552, 10).Value = 600
$ws.Cells.Item(552, 11).Value = 7500
$ws.Cells.Item(552, 12).Value = 7500
$ws.Cells.Item(552, 13).Value = 7500
$ws.Cells.Item(552, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(552, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(552, 16).Value = 375
$ws.Cells.Item(552, 17).Value = 20
$ws.Cells.Item(552, 18).Value = "Hortaliza"
